# journal de travail MA-20 - add new journal entries (score table / bug fixes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: fill in the previously-empty date (A20) ---
# Copy the date formatting (style s="4", numFmtId 14) from A13 so we reuse
# the existing cell style instead of minting a new one.
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A20").Value = 43908                     # 2020-03-18

# --- Row 21: new journal entry (was fully blank except format on E21) ---
# Values are entered in the same order the original author typed them so the
# shared-string table gets the same append order (F, then C, D, E).
$ws.Range("F21").Value = "Un problème de switch a été trouvé. Le bug consistait que si on choisis de quitter la partie quand on a gagné, le programme va rediriger sur la partie déjà existante mais qui est déjà finie. J'ai résolu le problème en remplaçant le switch par des if avec des break à l'intérieur"
$ws.Range("C21").Value = "Création du tableaux des scores, création du fichier contenant les noms d'utilisateurs ainsi que leurs scores et correction de bugs divers."
$ws.Range("D21").Value = "3h10"
$ws.Range("E21").Value = "Explication d'une découverte de bug juste à droite ainsi qu'un autre que lorsque on retire sur un bateau une fois ses points de vie descendent meme si la case a déjà été touchée. Le bug a été corrigé. J'ai aussi regardé une suite de vidéos explicant les diverses fonctions concernant la création de fichiers."

# A21 needs the same date style/format as the other date cells (s="4"),
# B21 already carries the right style, as do D21/E21.
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A21").Value = 43912                     # 2020-03-22
$ws.Range("B21").Value = 6

# C21 and F21 need the wrap-text "Problèmes/Résolution"-style format (s="5"),
# matching C20/E20/E21's existing style - copy it over (PasteSpecial with
# xlPasteFormats only touches formatting, the values set above are kept).
$ws.Range("E20").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E20").Copy() | Out-Null
$ws.Range("F21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- cosmetic view updates recorded in the diff ---
$excel.ActiveWindow.Zoom = 94
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E21").Select() | Out-Null
